$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added DB table columns
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "LastName"

# Row 2 - z5156156 / j.meraachli
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = "Jay"
$ws.Range("E2").Value = "Meraachli"

# Row 3 - z5020362 / k.schroder-turner
$ws.Range("C3").Value = 456
$ws.Range("D3").Value = "Khan"
$ws.Range("E3").Value = "Schroder"

# Row 4 - z5240067 / k.sharma.1
$ws.Range("C4").Value = 789
$ws.Range("D4").Value = "Kovid"
$ws.Range("E4").Value = "Sharma"

# Row 5 - z3264122 / simon
$ws.Range("C5").Value = 123
$ws.Range("D5").Value = "Simon"
$ws.Range("E5").Value = "Garrod"

# Row 6 - z5233368 / t.brunette
$ws.Range("C6").Value = 456
$ws.Range("D6").Value = "Tim"
$ws.Range("E6").Value = "Brunette"

# Widen column B to fit the longest email address, matching the author's resize
# (target stored width is 44.42578125 "characters"; the engine re-quantizes the
# value to its own pixel grid, so back the stored offset out of the request)
$ws.Columns.Item(2).ColumnWidth = 43.592447916666664

# Match the author's final cursor/selection position when they saved
[void]$ws.Range("G13").Select()
